# Add chemistry rules for `Fire` and `Ice` cards.
# Updates player roster rows 2-12 (Name, Rarity, Position, Country, League,
# Club, Cost, Chemistry, Org_Row_ID) to reflect the new card pool.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Mario Hermoso"
$ws.Range("E2").Value = "CB"
$ws.Range("G2").Value = "Spain"
$ws.Range("H2").Value = "LALIGA EA SPORTS"
$ws.Range("I2").Value = "Atlético de Madrid"
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 48

# Row 3
$ws.Range("A3").Value = "David García"
$ws.Range("D3").Value = "Common"
$ws.Range("E3").Value = "CB"
$ws.Range("G3").Value = "Spain"
$ws.Range("H3").Value = "LALIGA EA SPORTS"
$ws.Range("I3").Value = "CA Osasuna"
$ws.Range("N3").Value = 700
$ws.Range("P3").Value = 61

# Row 4
$ws.Range("A4").Value = "Mertens"
$ws.Range("B4").Value = 81
$ws.Range("E4").Value = "ST"
$ws.Range("G4").Value = "Belgium"
$ws.Range("H4").Value = "Trendyol Süper Lig"
$ws.Range("I4").Value = "Galatasaray"
$ws.Range("N4").Value = 700
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 108

# Row 5
$ws.Range("A5").Value = "Fred"
$ws.Range("E5").Value = "CDM"
$ws.Range("G5").Value = "Brazil"
$ws.Range("H5").Value = "Trendyol Süper Lig"
$ws.Range("I5").Value = "Fenerbahçe"
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 111

# Row 6
$ws.Range("A6").Value = "Lo Celso"
$ws.Range("E6").Value = "ST"
$ws.Range("G6").Value = "Argentina"
$ws.Range("H6").Value = "Premier League"
$ws.Range("I6").Value = "Spurs"
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 115

# Row 7
$ws.Range("A7").Value = "Kamara"
$ws.Range("E7").Value = "CDM"
$ws.Range("G7").Value = "France"
$ws.Range("H7").Value = "Premier League"
$ws.Range("I7").Value = "Aston Villa"
$ws.Range("P7").Value = 118

# Row 8
$ws.Range("A8").Value = "Lamela"
$ws.Range("B8").Value = 80
$ws.Range("E8").Value = "CAM"
$ws.Range("G8").Value = "Argentina"
$ws.Range("H8").Value = "LALIGA EA SPORTS"
$ws.Range("I8").Value = "Sevilla FC"
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = 125

# Row 9
$ws.Range("A9").Value = "Digne"
$ws.Range("E9").Value = "LB"
$ws.Range("H9").Value = "Premier League"
$ws.Range("I9").Value = "Aston Villa"
$ws.Range("N9").Value = 800
$ws.Range("O9").Value = 3
$ws.Range("P9").Value = 147

# Row 10
$ws.Range("A10").Value = "Nacho Vidal"
$ws.Range("B10").Value = 76
$ws.Range("E10").Value = "RB"
$ws.Range("G10").Value = "Spain"
$ws.Range("H10").Value = "LALIGA EA SPORTS"
$ws.Range("I10").Value = "CA Osasuna"
$ws.Range("N10").Value = 700
$ws.Range("O10").Value = 3
$ws.Range("P10").Value = 195

# Row 11
$ws.Range("A11").Value = "Lincoln"
$ws.Range("B11").Value = 75
$ws.Range("E11").Value = "CAM"
$ws.Range("G11").Value = "Brazil"
$ws.Range("H11").Value = "Trendyol Süper Lig"
$ws.Range("I11").Value = "Fenerbahçe"
$ws.Range("N11").Value = 650
$ws.Range("O11").Value = 3
$ws.Range("P11").Value = 223

# Row 12
$ws.Range("A12").Value = "Barrea"
$ws.Range("B12").Value = 63
$ws.Range("G12").Value = "Argentina"
$ws.Range("H12").Value = "LPF"
$ws.Range("I12").Value = "Godoy Cruz"
$ws.Range("P12").Value = 239
